# OrangeHRMS.xlsx edit — "adding code to check success msg"
#
# The workbook has two worksheets. Note the sheet *names* are swapped
# relative to their actual content/backing part:
#   - Worksheets.Item(1), named "Test Cases", is the small Runmode
#     master list (A1:B5) driving which top-level test cases run.
#   - Worksheets.Item(2), named "Test Data", is the big per-test-case
#     data table (A1:G31).
#
# The edit flips several Runmode flags from "Y" to "N" (skip these
# tests), and updates the saved cursor/selection state to match.

$wb = $excel.ActiveWorkbook

$wsRunmodeList = $wb.Worksheets.Item(1)   # "Test Cases" sheet / master Runmode list
$wsTestData    = $wb.Worksheets.Item(2)   # "Test Data" sheet / detailed test data

# --- Master Runmode list: stop running DeleteJobTitleTest & EditJobTitleTest ---
$wsRunmodeList.Range("B4").Value = "N"
$wsRunmodeList.Range("B5").Value = "N"

# --- JobTitleTest block (rows 9-14 of the G "Runmode" column) ---
$wsTestData.Range("G9").Value  = "N"
$wsTestData.Range("G10").Value = "N"
$wsTestData.Range("G11").Value = "N"
$wsTestData.Range("G12").Value = "N"
$wsTestData.Range("G13").Value = "N"
$wsTestData.Range("G14").Value = "N"

# --- DeleteTitleTest block (rows 20-21 of the D "Runmode" column) ---
$wsTestData.Range("D20").Value = "N"
$wsTestData.Range("D21").Value = "N"

# --- Restore cursor / selection state left by the editing session ---
$wsRunmodeList.Activate()
$wsRunmodeList.Range("B3").Select()

$wsTestData.Activate()
# Scroll the window so row 7 is the first visible row (best effort —
# persisted as topLeftCell when supported by the host).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsTestData.Range("D19").Select()
